$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly rows of Fruta/Hortaliza data got rotated: row2 <- old row3,
# row3 <- old row6, row6 <- old row8, row8 <- old row2 (cyclic shift of
# Fecha/Volumen/Precio fields, everything else in the row stays the same).

# Capture the original values before overwriting anything (use Value2 to
# read raw numeric values).
$orig2 = @{
    D = $ws.Cells.Item(2, 4).Value2
    M = $ws.Cells.Item(2, 13).Value2
    N = $ws.Cells.Item(2, 14).Value2
    O = $ws.Cells.Item(2, 15).Value2
    P = $ws.Cells.Item(2, 16).Value2
    S = $ws.Cells.Item(2, 19).Value2
}
$orig3 = @{
    D = $ws.Cells.Item(3, 4).Value2
    M = $ws.Cells.Item(3, 13).Value2
    N = $ws.Cells.Item(3, 14).Value2
    O = $ws.Cells.Item(3, 15).Value2
    P = $ws.Cells.Item(3, 16).Value2
    S = $ws.Cells.Item(3, 19).Value2
}
$orig6 = @{
    D = $ws.Cells.Item(6, 4).Value2
    M = $ws.Cells.Item(6, 13).Value2
    N = $ws.Cells.Item(6, 14).Value2
    O = $ws.Cells.Item(6, 15).Value2
    P = $ws.Cells.Item(6, 16).Value2
    S = $ws.Cells.Item(6, 19).Value2
}
$orig8 = @{
    D = $ws.Cells.Item(8, 4).Value2
    M = $ws.Cells.Item(8, 13).Value2
    N = $ws.Cells.Item(8, 14).Value2
    O = $ws.Cells.Item(8, 15).Value2
    P = $ws.Cells.Item(8, 16).Value2
    S = $ws.Cells.Item(8, 19).Value2
}

function Set-RowValues($rowNum, $vals) {
    $ws.Cells.Item($rowNum, 4).Value = $vals.D
    $ws.Cells.Item($rowNum, 13).Value = $vals.M
    $ws.Cells.Item($rowNum, 14).Value = $vals.N
    $ws.Cells.Item($rowNum, 15).Value = $vals.O
    $ws.Cells.Item($rowNum, 16).Value = $vals.P
    $ws.Cells.Item($rowNum, 19).Value = $vals.S
}

Set-RowValues 2 $orig3
Set-RowValues 3 $orig6
Set-RowValues 6 $orig8
Set-RowValues 8 $orig2
